# Dodat rad sa objektima
# Dodati modeli za objekat, rad sa bazom objekata, zapocet rad sa canvasom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "radnik" model header moves from col idx 25 -> reused string (col F) ---
$ws.Range("F1").Value = "radnik"

# --- Write brand-new shared strings first, in the exact order they must appear
#     in the shared string table (idR, Array[Objekat], idSkica, _id, vlasnik) ---
$ws.Range("F2").Value = "idR"
$ws.Range("B8").Value = "Array[Objekat]"
$ws.Range("E8").Value = "idSkica"
$ws.Range("E2").Value = "_id"
$ws.Range("E3").Value = "vlasnik"

# --- Remaining cells: existing strings reused, reflowed into their new positions ---
$ws.Range("F3").Value = "ime"
$ws.Range("E4").Value = "tip"
$ws.Range("F4").Value = "prezime"
$ws.Range("E5").Value = "adresa"
$ws.Range("F5").Value = "mejl"
$ws.Range("E6").Value = "brProstorija"
$ws.Range("F6").Value = "telefon"
$ws.Range("E7").Value = "kvadratura"
$ws.Range("F7").Value = "specijalizacija"
$ws.Range("C12").Value = "brojRadnika"

# --- Column B width: 13.7109375 -> 16 (character width units ~15.14 maps to XML width 16) ---
$ws.Columns.Item(2).ColumnWidth = 15.14

# --- Selection moves from C13 to E4 ---
$ws.Range("E4").Select()
